# CMG cutover template reset: clear out the previous cutover's captured
# metrics and reorder the SAEGW label rows so SAEGW35/36 are queued up
# next, per "modify excel for cmg add saegw35,36".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the captured per-host metric values (rows 3-12, columns B:H) while
# keeping the existing cell formatting/styles intact.
$ws.Range("B3:H12").ClearContents()

# Rows 16-18 list the three SAEGW hosts being tracked along with their
# sample command output values. Re-point the host names (35/36 now come
# first, pushing 37 to the top) and clear the stale sample values.
$ws.Range("A16").Value = "SHSAEGW37BNK"
$ws.Range("A17").Value = "SHSAEGW36BNK"
$ws.Range("A18").Value = "SHSAEGW35BNK"
$ws.Range("B16:F16").ClearContents()
$ws.Range("B17:F17").ClearContents()
$ws.Range("B18:F18").ClearContents()

# Move the selection/scroll position to A18 (clears the stale
# topLeftCell scroll anchor left over from the prior edit session).
$ws.Range("A18").Select()
